$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-looking-numeric values to be stored as text (matching the
# workbook's existing convention, e.g. A2/B2 = "123" as text), then
# restore the number format so no stray formatting is left behind.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "123456"
$ws.Range("A3").ClearFormats()

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "111"
$ws.Range("B3").ClearFormats()

$ws.Range("C3").Value = "Cliente"

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2"
$ws.Range("A4").ClearFormats()

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2"
$ws.Range("B4").ClearFormats()

$ws.Range("C4").Value = "Cliente"
